$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 0.93932292750621427
$ws.Range("G2").Value = 0.97080663208218898
$ws.Range("H2").Value = 0.96590909090909083

# Row 3
$ws.Range("E3").Value = 0.89040412662537027
$ws.Range("G3").Value = 0.96607909939268255
$ws.Range("H3").Value = 0.91304347826086951

# Row 4
$ws.Range("D4").Value = 0.82
$ws.Range("E4").Value = 0.94072222115925452
$ws.Range("G4").Value = 0.99453551912568305
$ws.Range("H4").Value = 0.94444444444444442

# Row 5
$ws.Range("D5").Value = 0.998
$ws.Range("E5").Value = 0.87934215774378044
$ws.Range("F5").Value = 0.98499999999999999
$ws.Range("G5").Value = 0.89285714285714279
$ws.Range("H5").Value = 0.88

# Row 7 (Mac example text + thresholds/results)
$ws.Range("C7").Value = "The example is about Mac"
$ws.Range("D7").Value = 0.98199999999999998
$ws.Range("E7").Value = 0.76394368355648257
$ws.Range("G7").Value = 0.91151202749140903
$ws.Range("H7").Value = 0.76923076923076916

# Row 8 (iCloud example text + thresholds/results)
$ws.Range("C8").Value = "The sentence is about icloud"
$ws.Range("D8").Value = 0.94
$ws.Range("E8").Value = 0.81443211092126233
$ws.Range("F8").Value = 0.995
$ws.Range("G8").Value = 0.83333333333333326
$ws.Range("H8").Value = 0.8

# Row 9 (iTunes example text + thresholds/results)
$ws.Range("C9").Value = "The example is iTunes"
$ws.Range("D9").Value = 0.98299999999999998
$ws.Range("E9").Value = 0.86382457327921347
$ws.Range("F9").Value = 0.995
$ws.Range("G9").Value = 0.875
$ws.Range("H9").Value = 0.8571428571428571
